$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2:B4 currently hold the text "2" (shared string). They need to hold
# the text "1" instead - still a text value (shared string), not a
# number. Plainly assigning .Value = "1" would be auto-coerced to the
# number 1 by Excel, so the cells are briefly marked as Text, the
# literal string is entered, and the formatting is cleared again so the
# cells end up back at the default (General) style while keeping "1"
# stored as text.
$rng = $ws.Range("B2:B4")
$rng.NumberFormat = "@"
$rng.Value = "1"
$rng.ClearFormats()
